$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.282872386776422
$ws.Range("C2").Value = 0.1391894435007327
$ws.Range("D2").Value = -0.182694870696341
$ws.Range("E2").Value = 0.3443354285759472
$ws.Range("F2").Value = 0.1070053646347159
$ws.Range("G2").Value = 0.2095078597799856
$ws.Range("H2").Value = 0.110633009540361
$ws.Range("I2").Value = 0.6128907522960124
$ws.Range("J2").Value = 0.8708754778845325
$ws.Range("K2").Value = -0.8533573809972426
$ws.Range("B3").Value = 0.2750606270018259
$ws.Range("C3").Value = 0.5534868575902183
$ws.Range("D3").Value = 0.2080306981463169
$ws.Range("E3").Value = 0.2637694639246414
$ws.Range("F3").Value = 0.1446758428348626
$ws.Range("G3").Value = 0.6381972808990434
$ws.Range("H3").Value = 0.892409660067914
$ws.Range("I3").Value = -0.8334509197017121
$ws.Range("J3").Value = -0.1806300843769859
$ws.Range("K3").Value = 0.1825572214681606
$ws.Range("B4").Value = 0.1918994249285815
$ws.Range("C4").Value = 0.239828553585737
$ws.Range("D4").Value = 0.1203502894047785
$ws.Range("E4").Value = 0.6162427231036979
$ws.Range("F4").Value = 0.8725296081739512
$ws.Range("G4").Value = -0.8520779198862551
$ws.Range("H4").Value = -0.1986054693098328
$ws.Range("I4").Value = 0.164895676365706
$ws.Range("J4").Value = -0.7358907623030519
$ws.Range("K4").Value = 0.3344822619593895
$ws.Range("B5").Value = 0.2023032603298698
$ws.Range("C5").Value = 0.6254400633356464
$ws.Range("D5").Value = 0.85282204567298
$ws.Range("E5").Value = -0.8813162353832781
$ws.Range("F5").Value = -0.2307498541254956
$ws.Range("G5").Value = 0.1318997099156809
$ws.Range("H5").Value = -0.7691308072183836
$ws.Range("I5").Value = 0.301173166579823
$ws.Range("J5").Value = -0.2759067096457259
$ws.Range("K5").Value = -0.4815182333746927
$ws.Range("B6").Value = 0.8240207549100882
$ws.Range("C6").Value = -0.8834400190618514
$ws.Range("D6").Value = -0.2234668893430731
$ws.Range("E6").Value = 0.1423162228589251
$ws.Range("F6").Value = -0.7576896617840617
$ws.Range("G6").Value = 0.3129479666538719
$ws.Range("H6").Value = -0.2640230888223336
$ws.Range("I6").Value = -0.4695989885331262
$ws.Range("J6").Value = 0.5095561729841656
$ws.Range("K6").Value = -0.226786726508225
$ws.Range("B7").Value = -0.3453490504388863
$ws.Range("C7").Value = 0.05047105392560119
$ws.Range("D7").Value = -0.83756240835547
$ws.Range("E7").Value = 0.2380186391052039
$ws.Range("F7").Value = -0.3369057261414908
$ws.Range("G7").Value = -0.5416351540011936
$ws.Range("H7").Value = 0.4378696407944773
$ws.Range("I7").Value = -0.2983290036150871
$ws.Range("J7").Value = -0.3162381895859678
$ws.Range("K7").Value = -0.2928084366771486
$ws.Range("B8").Value = -0.7422476366981072
$ws.Range("C8").Value = 0.3162863437924887
$ws.Range("D8").Value = -0.2663152355556169
$ws.Range("E8").Value = -0.4745738317060451
$ws.Range("F8").Value = 0.5033030119909454
$ws.Range("G8").Value = -0.2336496331081919
$ws.Range("H8").Value = -0.2519093219356206
$ws.Range("I8").Value = -0.2286430458875994
$ws.Range("J8").Value = -0.5938837416006923
$ws.Range("K8").Value = 0.06190409654164158
$ws.Range("B9").Value = -0.00619269375287701
$ws.Range("C9").Value = -0.3336679001236204
$ws.Range("D9").Value = 0.5887475516261532
$ws.Range("E9").Value = -0.1739978958108406
$ws.Range("F9").Value = -0.2042525916245146
$ws.Range("G9").Value = -0.1865645222487989
$ws.Range("H9").Value = -0.5543992815454399
$ws.Range("I9").Value = 0.1001822482843727
$ws.Range("J9").Value = -0.03126473275901087
$ws.Range("K9").Value = -0.1410506211186619
$ws.Range("B10").Value = 0.3911393015443105
$ws.Range("C10").Value = -0.2979577194736586
$ws.Range("D10").Value = -0.2956322122566587
$ws.Range("E10").Value = -0.2635227406874879
$ws.Range("F10").Value = -0.6249641526670209
$ws.Range("G10").Value = 0.03245628485901808
$ws.Range("H10").Value = -0.09772803603944658
$ws.Range("I10").Value = -0.2069514014028122
$ws.Range("J10").Value = -1.036602919657539
$ws.Range("K10").Value = -0.4560229796881132
$ws.Range("B11").Value = -0.228230158963832
$ws.Range("C11").Value = -0.1886450477650728
$ws.Range("D11").Value = -0.5467980559473664
$ws.Range("E11").Value = 0.1120680736150478
$ws.Range("F11").Value = -0.01748087635905388
$ws.Range("G11").Value = -0.1264250868990488
$ws.Range("H11").Value = -0.9559539928657461
$ws.Range("I11").Value = -0.3753202137097185
$ws.Range("J11").Value = -0.4586594580766912
$ws.Range("K11").Value = -0.6765490042469573
$ws.Range("B12").Value = -0.5591547568613399
$ws.Range("C12").Value = 0.09787708368923709
$ws.Range("D12").Value = -0.03242713058169561
$ws.Range("E12").Value = -0.1416815451927503
$ws.Range("F12").Value = -0.9713369016796535
$ws.Range("G12").Value = -0.390754202413632
$ws.Range("H12").Value = -0.4741138500012969
$ws.Range("I12").Value = -0.6920114308521601
$ws.Range("J12").Value = 0.8717720887233864
$ws.Range("K12").Value = -0.4642089603862817
$ws.Range("B13").Value = 0.1398755720077088
$ws.Range("C13").Value = -0.04718941243373825
$ws.Range("D13").Value = -0.9130114627925084
$ws.Range("E13").Value = -0.3492324732466185
$ws.Range("F13").Value = -0.4403994753588329
$ws.Range("G13").Value = -0.6619245159313032
$ws.Range("H13").Value = 0.900173610209962
$ws.Range("I13").Value = -0.4365905078372946
$ws.Range("J13").Value = 1.003661313589152
$ws.Range("K13").Value = -0.09695126965353879
$ws.Range("B14").Value = -1.007520009675017
$ws.Range("C14").Value = -0.3857468722612991
$ws.Range("D14").Value = -0.4502017751758141
$ws.Range("E14").Value = -0.6594324126472006
$ws.Range("F14").Value = 0.9083245926860781
$ws.Range("G14").Value = -0.4258347558073555
$ws.Range("H14").Value = 1.015616079910623
$ws.Range("I14").Value = -0.08444455851082849
$ws.Range("J14").Value = -0.2945431432588089
$ws.Range("K14").Value = 0.3040816658791113
$ws.Range("B15").Value = -0.3077970529158559
$ws.Range("C15").Value = -0.5788562590806742
$ws.Range("D15").Value = 0.961338929688949
$ws.Range("E15").Value = -0.3851068719675381
$ws.Range("F15").Value = 1.050866189894518
$ws.Range("G15").Value = -0.05163699300024938
$ws.Range("H15").Value = -0.2628248666403167
$ws.Range("I15").Value = 0.3353140868060093
$ws.Range("J15").Value = 0.009465662865194846
$ws.Range("K15").Value = 0.7265657060874963
$ws.Range("B16").Value = 1.139271322614915
$ws.Range("C16").Value = -0.2982492763927073
$ws.Range("D16").Value = 1.093791915679716
$ws.Range("E16").Value = -0.02972831832925377
$ws.Range("F16").Value = -0.250970648745674
$ws.Range("G16").Value = 0.3423607205295725
$ws.Range("H16").Value = 0.01421452142035373
$ws.Range("I16").Value = 0.7302167878130713
$ws.Range("J16").Value = 2.425529404116681
$ws.Range("K16").Value = 9.269510203911928
$ws.Range("B17").Value = -0.2876812562600493
$ws.Range("C17").Value = 1.102566704251488
$ws.Range("D17").Value = -0.02192591428551849
$ws.Range("E17").Value = -0.2436605879111656
$ws.Range("F17").Value = 0.3494215681832005
$ws.Range("G17").Value = 0.02114969792139543
$ws.Range("H17").Value = 0.7370887910462527
$ws.Range("I17").Value = 2.432369740824074
$ws.Range("J17").Value = 9.276334707790259
$ws.Range("K17").Value = -8.274453695494744
$ws.Range("B18").Value = 0.9955923004358076
$ws.Range("C18").Value = -0.0871855372237873
$ws.Range("D18").Value = -0.2892511864441449
$ws.Range("E18").Value = 0.3130908916669468
$ws.Range("F18").Value = -0.0108203968125144
$ws.Range("G18").Value = 0.7071725634927386
$ws.Range("H18").Value = 2.403421102624959
$ws.Range("I18").Value = 9.247842001542866
$ws.Range("J18").Value = -8.302731520459806
$ws.Range("K18").Value = -0.5932976437114486
$ws.Range("B19").Value = -0.0533123394792912
$ws.Range("C19").Value = -0.2864980261873478
$ws.Range("D19").Value = 0.3015019608161922
$ws.Range("E19").Value = -0.02894181086942849
$ws.Range("F19").Value = 0.6860724501299835
$ws.Range("G19").Value = 2.380962765938485
$ws.Range("H19").Value = 9.224764288378545
$ws.Range("I19").Value = -8.326091704608711
$ws.Range("J19").Value = -0.6167866614867294
$ws.Range("K19").Value = 0.9660844828649828
$ws.Range("B20").Value = -0.5437734556935978
$ws.Range("C20").Value = 0.1702226893339956
$ws.Range("D20").Value = -0.1022974921431762
$ws.Range("E20").Value = 0.63883117546281
$ws.Range("F20").Value = 2.345515881559747
$ws.Range("G20").Value = 9.194643401784708
$ws.Range("H20").Value = -8.35380748862633
$ws.Range("I20").Value = -0.6434163542466904
$ws.Range("J20").Value = 0.9399452457891766
$ws.Range("K20").Value = -2.026139092643445
$ws.Range("B21").Value = 0.1269585753652727
$ws.Range("C21").Value = -0.1306388881700511
$ws.Range("D21").Value = 0.6141551499506877
$ws.Range("E21").Value = 2.321229248582915
$ws.Range("F21").Value = 9.170020105529264
$ws.Range("G21").Value = -8.378790661567026
$ws.Range("H21").Value = -0.668646313777596
$ws.Range("I21").Value = 0.9145696739774734
$ws.Range("J21").Value = -2.051594272676513
$ws.Range("K21").Value = 0.2124867959412257
$ws.Range("B22").Value = -0.02044773294388247
$ws.Range("C22").Value = 0.6702777838950345
$ws.Range("D22").Value = 2.355311064516826
$ws.Range("E22").Value = 9.195119489401716
$ws.Range("F22").Value = -8.357352010851283
$ws.Range("G22").Value = -0.6486996043872704
$ws.Range("H22").Value = 0.933908325411367
$ws.Range("I22").Value = -2.032503447757091
$ws.Range("J22").Value = 0.2314766118486793
$ws.Range("K22").Value = -0.2196984654044191
$ws.Range("B23").Value = 0.5134276986299353
$ws.Range("C23").Value = 2.260143333751186
$ws.Range("D23").Value = 9.140692358735786
$ws.Range("E23").Value = -8.393905998263049
$ws.Range("F23").Value = -0.6763647928704721
$ws.Range("G23").Value = 0.9105011857394278
$ws.Range("H23").Value = -2.053840367633118
$ws.Range("I23").Value = 0.2111440141036455
$ws.Range("J23").Value = -0.2395424444363658
$ws.Range("K23").Value = -0.379677132609091
$ws.Range("B24").Value = 2.211279720138829
$ws.Range("C24").Value = 9.112923508459115
$ws.Range("D24").Value = -8.411125843300541
$ws.Range("E24").Value = -0.688459908824067
$ws.Range("F24").Value = 0.9009170582553292
$ws.Range("G24").Value = -2.06219729845007
$ws.Range("H24").Value = 0.2033872967055401
$ws.Range("I24").Value = -0.2470056675122654
$ws.Range("J24").Value = -0.386996833524363
$ws.Range("K24").Value = 0.1394248687261353
